# QA test cases.xlsx - "Did final run of all tests and passed"
# Fill in the "Outcome" column (F) for the test rows that were run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Test Scenario: To Do App (Navigation and Login) -- rows 7-11
$ws.Range("F7").Value  = "User cannot login, error message displayed"
$ws.Range("F8").Value  = "User cannot login, error message displayed"
$ws.Range("F9").Value  = "User can log out, message displayed"
$ws.Range("F10").Value = "user cannot navigate, error message displayed"
$ws.Range("F11").Value = "user cannot navigate, error message displayed"

# Test Scenario: To Do Page -- rows 15-19
$ws.Range("F15").Value = "Displayed"
$ws.Range("F16").Value = "Item not added, error message shown"
$ws.Range("F17").Value = "Item added"
$ws.Range("F18").Value = "Item added"
$ws.Range("F19").Value = "Item deleted and removed from list"

# Last edited cell becomes the active selection
$ws.Range("F11").Select()
